$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.137.68'
$ws.Range("E2").Value = '  +10.89%  '
$ws.Range("D3").Value = '3.351.62'
$ws.Range("E3").Value = '  +5.28%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.33'
$ws.Range("E5").Value = '  +4.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '643.91'
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.332'
$ws.Range("E7").Value = '  +20.34%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.618'
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("D10").Value = '3.351.70'
$ws.Range("E10").Value = '  +5.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.599'
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000277'
$ws.Range("E12").Value = '  +7.11%  '
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.32'
$ws.Range("E14").Value = '  +10.79%  '
$ws.Range("D15").Value = '3.963.48'
$ws.Range("E15").Value = '  +5.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.42'
$ws.Range("E16").Value = '  +3.31%  '
$ws.Range("D17").Value = '88.420.91'
$ws.Range("E17").Value = '  +10.16%  '
$ws.Range("D18").Value = '3.334.51'
$ws.Range("E18").Value = '  +5.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.68'
$ws.Range("E19").Value = '  +3.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.22'
$ws.Range("E20").Value = '  +6.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '471.87'
$ws.Range("E21").Value = '  +7.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.15'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.44'
$ws.Range("E23").Value = '  +4.63%  '
$ws.Range("B24").Value = 'Aptos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.50'
$ws.Range("E24").Value = '  +24.16%  '
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.43'
$ws.Range("E25").Value = '  +7.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.39'
$ws.Range("E26").Value = '  +15.30%  '
$ws.Range("D27").Value = '3.507.90'
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '79.31'
$ws.Range("E28").Value = '  +4.45%  '
$ws.Range("E29").Value = '  +77.40%  '
$ws.Range("E30").Value = '  +6.75%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '599.82'
$ws.Range("E32").Value = '  +7.53%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.32'
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").Value = '  +7.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.07'
$ws.Range("E36").Value = '  +3.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.151'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.95'
$ws.Range("E38").Value = '  +4.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.79'
$ws.Range("E39").Value = '  +20.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.422'
$ws.Range("E40").Value = '  +3.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.996'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.46'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.10'
$ws.Range("E43").Value = '  +16.25%  '
$ws.Range("E44").Value = '  +13.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.24'
$ws.Range("E45").Value = '  -2.93%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '192.36'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.37'
$ws.Range("E48").Value = '  +10.41%  '
$ws.Range("E49").Value = '  +5.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.798'
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.71'
$ws.Range("E51").Value = '  +6.99%  '
